$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new value. Values in the Price (D) column
# that look like plain numbers are prefixed with a leading apostrophe so
# Excel stores them as text (matching the source data which uses
# dotted/European-style numbers in neighbouring cells).
$data = @{
    2 = @{ D="65.188.54"; E="  -0.47%  " }
    3 = @{ D="3.552.56"; E="  -0.19%  " }
    4 = @{ E="  -0.03%  " }
    5 = @{ D="'597.87"; E="  -0.09%  " }
    6 = @{ D="'133.49"; E="  -5.34%  " }
    7 = @{ D="3.552.72"; E="  -0.18%  " }
    8 = @{ E="  -0.01%  " }
    9 = @{ E="  -0.60%  " }
    10 = @{ E="  -2.53%  " }
    11 = @{ D="'7.11"; E="  -0.77%  " }
    12 = @{ E="  -0.98%  " }
    13 = @{ D="4.152.47"; E="  -0.23%  " }
    14 = @{ E="  -2.62%  " }
    15 = @{ D="'26.94"; E="  -0.56%  " }
    16 = @{ D="3.549.20"; E="  -0.13%  " }
    17 = @{ E="  -0.14%  " }
    18 = @{ D="65.274.92"; E="  -0.16%  " }
    19 = @{ D="'9.95"; E="  -4.66%  " }
    20 = @{ D="'14.36"; E="  +1.00%  " }
    21 = @{ E="  -0.91%  " }
    22 = @{ D="'391.40"; E="  -1.44%  " }
    23 = @{ D="'0.577"; E="  +1.21%  " }
    24 = @{ D="3.694.60"; E="  -0.13%  " }
    25 = @{ D="'74.04"; E="  -0.73%  " }
    26 = @{ E="  +0.01%  " }
    27 = @{ E="  -0.92%  " }
    28 = @{ D="'7.81"; E="  +0.56%  " }
    29 = @{ D="'1.56"; E="  +25.44%  " }
    30 = @{ E="  +0.33%  " }
    31 = @{ D="'8.55"; E="  +3.17%  " }
    32 = @{ E="  +0.04%  " }
    33 = @{ D="3.552.90"; E="  -0.57%  " }
    34 = @{ D="'24.08"; E="  +0.16%  " }
    35 = @{ E="  -0.01%  " }
    36 = @{ D="'0.147"; E="  -0.29%  " }
    37 = @{ D="'170.40"; E="  +1.04%  " }
    38 = @{ D="'6.94"; E="  -1.44%  " }
    39 = @{ E="  -0.41%  " }
    40 = @{ D="'5.03"; E="  +1.03%  " }
    41 = @{ D="'0.0811"; E="  +0.70%  " }
    42 = @{ D="'0.827"; E="  -0.05%  " }
    43 = @{ D="'26.55"; E="  +0.22%  " }
    44 = @{ B="ONDO"; C="https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"; D="'1.25"; E="  +5.19%  " }
    45 = @{ B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="'43.09"; E="  +0.99%  " }
    46 = @{ E="  -0.06%  " }
    47 = @{ E="  -0.14%  " }
    48 = @{ E="  -1.84%  " }
    49 = @{ D="2.457.05"; E="  +2.36%  " }
    50 = @{ D="'6.91"; E="  +1.21%  " }
    51 = @{ E="  +0.62%  " }
}

$colIndex = @{ B = 2; C = 3; D = 4; E = 5 }

foreach ($r in $data.Keys) {
    $rowData = $data[$r]
    foreach ($col in $rowData.Keys) {
        $c = $colIndex[$col]
        $ws.Cells.Item([int]$r, $c).Value = $rowData[$col]
    }
}

Write-Host "Applied $($data.Count) row updates"
